# Applies cryptocurrency price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '25.944.18'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = "'" + '1.637.82'
$ws.Range("E3").Value = '  +0.05%  '

$ws.Range("E4").Value = '  +0.88%  '

$ws.Range("D5").Value = "'" + '214.75'
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("E6").Value = '  +0.68%  '

$ws.Range("E7").Value = '  +0.94%  '

$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("E9").Value = '  +0.62%  '

$ws.Range("D10").Value = "'" + '19.61'
$ws.Range("E10").Value = '  -0.93%  '

$ws.Range("D11").Value = "'" + '0.0794'
$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").Value = "'" + '1.864.44'
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'" + '1.663.18'
$ws.Range("E13").Value = '  +1.84%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = "'" + '4.25'
$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("E15").Value = '  -1.49%  '

$ws.Range("E16").Value = '  -0.86%  '

$ws.Range("D17").Value = "'" + '62.51'
$ws.Range("E17").Value = '  -1.06%  '

$ws.Range("D18").Value = "'" + '25.952.76'
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").Value = "'" + '193.55'
$ws.Range("E20").Value = '  +0.43%  '

$ws.Range("E21").Value = '  -1.56%  '

$ws.Range("E22").Value = '  -0.78%  '

$ws.Range("E23").Value = '  -1.55%  '

$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = "'" + '144.18'
$ws.Range("E24").Value = '  +1.82%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = "'" + '1.79'
$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("E26").Value = '  +0.98%  '

$ws.Range("E27").Value = '  +2.55%  '

$ws.Range("E28").Value = '  -0.54%  '

$ws.Range("D29").Value = "'" + '15.46'
$ws.Range("E29").Value = '  -0.90%  '

$ws.Range("E30").Value = '  -0.43%  '

$ws.Range("D31").Value = "'" + '0.0501'
$ws.Range("E31").Value = '  +1.52%  '

$ws.Range("E32").Value = '  -1.35%  '

$ws.Range("D33").Value = "'" + '3.23'
$ws.Range("E33").Value = '  -0.48%  '

$ws.Range("E34").Value = '  -2.75%  '

$ws.Range("E35").Value = '  +2.13%  '

$ws.Range("D36").Value = "'" + '0.902'
$ws.Range("E36").Value = '  -0.48%  '

$ws.Range("D37").Value = "'" + '1.138.18'
$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("D38").Value = "'" + '0.546'
$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("E39").Value = '  -1.29%  '

$ws.Range("E40").Value = '  +0.79%  '

$ws.Range("D41").Value = "'" + '99.42'
$ws.Range("E41").Value = '  -0.77%  '

$ws.Range("D42").Value = "'" + '0.799'
$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("D43").Value = "'" + '5.41'
$ws.Range("E43").Value = '  -2.58%  '

$ws.Range("D44").Value = "'" + '1.773.80'
$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("E45").Value = '  +7.74%  '

$ws.Range("D46").Value = "'" + '56.44'
$ws.Range("E46").Value = '  +1.45%  '

$ws.Range("E47").Value = '  +2.76%  '

$ws.Range("E48").Value = '  +0.37%  '

$ws.Range("D49").Value = "'" + '7.66'
$ws.Range("E49").Value = '  +0.33%  '

$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("D51").Value = "'" + '0.0962'
